# Adds a new "رقم_المعاملة_المرجعية" (reference transaction number) column (K)
# to the transactions sheet, and appends 3 new transaction rows (202-204)
# that record a correction: the original entry (202) is reversed via an
# "خروج" row (203) and replaced with a corrected "دخول" row (204), both of
# which point back at the original transaction through the new K column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header (K1), formatted like the rest of the header row ---
$ws.Range("K1").Value = "رقم_المعاملة_المرجعية"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 202: original (erroneous) stock-in entry ---
$ws.Range("A202").Value = 1765132405
$ws.Range("B202").Value = "مخزن_المستلزمات_الطبية"
$ws.Range("C202").Value = "2025-12-07 20:33:25"
$ws.Range("D202").Value = "طلاء أبيض"
$ws.Range("E202").Value = "مواد التشطيب"
$ws.Range("F202").Value = "دخول"
$ws.Range("G202").Value = 22
$ws.Range("I202").Value = 365

# --- Row 203: system-generated reversal ("خروج") of the row-202 entry ---
$ws.Range("A203").Value = 1765132427
$ws.Range("B203").Value = "مخزن_المستلزمات_الطبية"
$ws.Range("C203").Value = "2025-12-07 20:33:25"
$ws.Range("D203").Value = "طلاء أبيض"
$ws.Range("E203").Value = "مواد التشطيب"
$ws.Range("F203").Value = "خروج"
$ws.Range("G203").Value = 22
$ws.Range("H203").Value = "النظام"
$ws.Range("J203").Value = "إلغاء معاملة رقم 1765132405 - تصحيح خطأ في الإدخال"
$ws.Range("K203").Value = 1765132405

# --- Row 204: system-generated corrected ("دخول") replacement entry ---
$ws.Range("A204").Value = 1765132427
$ws.Range("B204").Value = "مخزن_المستلزمات_الطبية"
$ws.Range("C204").Value = "2025-12-07 20:33:25"
$ws.Range("D204").Value = "طلاء أبيض"
$ws.Range("E204").Value = "مواد التشطيب"
$ws.Range("F204").Value = "دخول"
$ws.Range("G204").Value = 15
$ws.Range("H204").Value = "النظام"
$ws.Range("J204").Value = "معاملة محدثة - تصحيح خطأ في الإدخال (من 22 إلى 15)"
$ws.Range("K204").Value = 1765132405
